$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bugs")

# Bug #2 (row 5) severity changes from P2 to P3
$ws.Range("D5").Value = "P3"

# New bug #3 (row 6)
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 45578
$ws.Range("B6").NumberFormat = "m/d/yy"
$ws.Range("C6").Value = "Yuntian"
$ws.Range("D6").Value = "P2"
$ws.Range("E6").Value = "resolved"
$ws.Range("F6").Value = "Yuntian"
$ws.Range("H6").Value = "The player should be blocked by the wall"
$ws.Range("I6").Value = "The player keeps moving"
$ws.Range("J6").Value = "If the player moves towards the bottom left corner of the map, they will pass through the wall"
$ws.Range("G6").Value = "The player can squeeze themselves into a specific corner and pass through the wall"

$ws.Rows.Item(6).RowHeight = 51

$ws.Range("G10").Select() | Out-Null
